$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "Main Category" -> "Category"
$ws.Range("A1").Value = "Category"

# Update "Poultry" row (row 11) values
$ws.Range("B11").Value = 19.41290322580645
$ws.Range("C11").Value = 9.129032258064516
$ws.Range("D11").Value = 2.474193548387097
$ws.Range("E11").Value = 68.74193548387096
$ws.Range("F11").Value = 204.5322580645161
$ws.Range("G11").Value = 0.4
$ws.Range("H11").Value = 170.1774193548387
$ws.Range("I11").Value = 92.34482758620689
$ws.Range("J11").Value = 31.25806451612903

# Insert a new row at 12 (pushes Starchy vegetables / Sweets and snacks / Vegetables down by one)
$ws.Rows.Item(12).Insert()

# Fill in the new "Red meat" row (row 12)
$ws.Range("A12").Value = "Red meat"
$ws.Range("B12").Value = 18.06887755102041
$ws.Range("C12").Value = 12.56428571428571
$ws.Range("D12").Value = 2.978571428571428
$ws.Range("E12").Value = 66.02551020408163
$ws.Range("F12").Value = 194.4974358974359
$ws.Range("G12").Value = 0.7729591836734694
$ws.Range("H12").Value = 198.0816326530612
$ws.Range("I12").Value = 98.80612244897959
$ws.Range("J12").Value = 33.97448979591837

# Match formatting/style of the other category column A cells (bold, bordered, centered)
$ws.Range("A13").Copy()
$ws.Range("A12").PasteSpecial(-4122)
